$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: Insert two new columns at D:E, shifting existing D:K to F:M
$ws.Range("D:E").Insert()

# Step 2: Fix formatting of the newly inserted D:E columns (inherit from F:G)
$ws.Range('F:G').Copy()
$ws.Range('D:E').PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 3: Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 528600
$ws.Range("E8").Value = 526800
$ws.Range("D9").Value = 356200
$ws.Range("E9").Value = 351700
$ws.Range("D10").Value = 172400
$ws.Range("E10").Value = 175100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 11400
$ws.Range("E15").Value = 11300
$ws.Range("D17").Value = 478500
$ws.Range("E17").Value = 484200
$ws.Range("D18").Value = 50100
$ws.Range("E18").Value = 42600
$ws.Range("D20").Value = 300
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 61800
$ws.Range("E21").Value = 53900
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = 4600
$ws.Range("D23").Value = 50400
$ws.Range("E23").Value = 38000
$ws.Range("D24").Value = 14800
$ws.Range("E24").Value = 10100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 35600
$ws.Range("E26").Value = 27900
$ws.Range("D27").Value = 35600
$ws.Range("E27").Value = 27900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -300
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 35600
$ws.Range("E33").Value = 27900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 35600
$ws.Range("E35").Value = 27900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 13900
$ws.Range("E41").Value = 18600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 416800
$ws.Range("E43").Value = 414200
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 51500
$ws.Range("E45").Value = 47000
$ws.Range("D46").Value = 482200
$ws.Range("E46").Value = 479800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 90400
$ws.Range("E48").Value = 86800
$ws.Range("D49").Value = 764700
$ws.Range("E49").Value = 771100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 155500
$ws.Range("E52").Value = 152700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1492700
$ws.Range("E54").Value = 1490400
$ws.Range("D57").Value = 31500
$ws.Range("E57").Value = 21200
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 275700
$ws.Range("E59").Value = 281900
$ws.Range("D60").Value = 307300
$ws.Range("E60").Value = 303100
$ws.Range("D61").Value = 440600
$ws.Range("E61").Value = 470400
$ws.Range("D62").Value = 105900
$ws.Range("E62").Value = 102200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 853700
$ws.Range("E66").Value = 875700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 286100
$ws.Range("E72").Value = 250400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 639000
$ws.Range("E76").Value = 614700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 35600
$ws.Range("E81").Value = 27900
$ws.Range("D83").Value = 11400
$ws.Range("E83").Value = 11300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 58900
$ws.Range("E89").Value = 42100
$ws.Range("D91").Value = -11300
$ws.Range("E91").Value = -7500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -7700
$ws.Range("E94").Value = -32700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -44300
$ws.Range("E100").Value = -36900
$ws.Range("D101").Value = 100
$ws.Range("E101").Value = 100
$ws.Range("D102").Value = 7100
$ws.Range("E102").Value = -27300

# Step 4: Apply data corrections identified in the shifted F:M range
$ws.Range("F61").Value = 475200
$ws.Range("F62").Value = 98100
$ws.Range("H89").Value = 58400
$ws.Range("I89").Value = 31300
$ws.Range("F91").Value = -10800
$ws.Range("J91").Value = -6200
$ws.Range("H102").Value = 33800
$ws.Range("I102").Value = 1900
